$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New enum "SkillId" added as row 10, with its items (rows 10-13) in the
# *items sub-table (columns G/I), mirroring the existing "AttributeType"
# enum pattern (row 4 header + rows 4-9 items).

# Row 10: enum header for SkillId (full_name, flags, unique) + first item (NONE = 0)
$ws.Range("B10").Value = "SkillId"
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = $true
$ws.Range("G10").Value = "NONE"
$ws.Range("I10").Value = 0

# Row 11: item CHAIN_LIGHTNING = 1
$ws.Range("G11").Value = "CHAIN_LIGHTNING"
$ws.Range("I11").Value = 1

# Row 12: item EXPLOSION = 2
$ws.Range("G12").Value = "EXPLOSION"
$ws.Range("I12").Value = 2

# Row 13: item SLOW = 3
$ws.Range("G13").Value = "SLOW"
$ws.Range("I13").Value = 3

# Update selection to match the final cursor position left after editing
[void]$ws.Range("K11").Select()
